# This script applies a permutation of data (columns D, H, I, J, K, L, M, N, O, P)
# across rows 3-18 and 20-22 of the active worksheet, matching the target diff.
# Row 19 and row 2 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each target row, the values to write (taken from the "before" state of
# the row indicated by the commit diff).
$rows = @{
    3  = @{ D = 44497; H = 'Sin especificar';        I = 'Primera'; J = 250; K = 800;  L = 800;  M = 800;  N = '$/kilo (volumen en unidades)'; O = 'Perú' ;                   P = 800  }
    4  = @{ D = 44510; H = 'Sin especificar';        I = 'Primera'; J = 250; K = 800;  L = 800;  M = 800;  N = '$/kilo (volumen en unidades)'; O = 'Perú' ;                   P = 800  }
    5  = @{ D = 44167; H = 'Sin especificar';        I = 'Primera'; J = 400; K = 5000; L = 5000; M = 5000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 5000 }
    6  = @{ D = 44167; H = 'Sin especificar';        I = 'Segunda'; J = 560; K = 3000; L = 3000; M = 3000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 3000 }
    7  = @{ D = 44167; H = 'Sin especificar';        I = 'Tercera'; J = 450; K = 2000; L = 2000; M = 2000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 2000 }
    8  = @{ D = 44495; H = 'Sin especificar';        I = 'Primera'; J = 200; K = 800;  L = 800;  M = 800;  N = '$/kilo (volumen en unidades)'; O = 'Perú' ;                   P = 800  }
    9  = @{ D = 44312; H = 'Sin especificar';        I = 'Primera'; J = 180; K = 2500; L = 2500; M = 2500; N = '$/unidad';                     O = 'Perú' ;                   P = 2500 }
    10 = @{ D = 44305; H = 'Sin especificar';        I = 'Primera'; J = 100; K = 2500; L = 2500; M = 2500; N = '$/unidad';                     O = 'Perú' ;                   P = 2500 }
    11 = @{ D = 44223; H = 'Americana O Klondike';   I = 'Extra';   J = 340; K = 2500; L = 2500; M = 2500; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 2500 }
    12 = @{ D = 44223; H = 'Americana O Klondike';   I = 'Primera'; J = 400; K = 2000; L = 2000; M = 2000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 2000 }
    13 = @{ D = 44223; H = 'Americana O Klondike';   I = 'Segunda'; J = 300; K = 1500; L = 1500; M = 1500; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 1500 }
    14 = @{ D = 44223; H = 'Americana O Klondike';   I = 'Tercera'; J = 160; K = 1000; L = 1000; M = 1000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 1000 }
    15 = @{ D = 44483; H = 'Sin especificar';        I = 'Primera'; J = 120; K = 800;  L = 800;  M = 800;  N = '$/kilo (volumen en unidades)'; O = 'Perú' ;                   P = 800  }
    16 = @{ D = 44488; H = 'Sin especificar';        I = 'Primera'; J = 150; K = 800;  L = 800;  M = 800;  N = '$/kilo (volumen en unidades)'; O = 'Perú' ;                   P = 800  }
    17 = @{ D = 44217; H = 'Sin especificar';        I = 'Extra';   J = 400; K = 2500; L = 2500; M = 2500; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 2500 }
    18 = @{ D = 44217; H = 'Sin especificar';        I = 'Primera'; J = 280; K = 2000; L = 2000; M = 2000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 2000 }
    20 = @{ D = 44504; H = 'Sin especificar';        I = 'Primera'; J = 200; K = 800;  L = 800;  M = 800;  N = '$/kilo (volumen en unidades)'; O = 'Perú' ;                   P = 800  }
    21 = @{ D = 44194; H = 'Sin especificar';        I = 'Extra';   J = 120; K = 3500; L = 3500; M = 3500; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 3500 }
    22 = @{ D = 44194; H = 'Sin especificar';        I = 'Primera'; J = 200; K = 3000; L = 3000; M = 3000; N = '$/unidad';                     O = "Región de O'Higgins" ;    P = 3000 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 8).Value  = $vals.H   # H - Variedad
    $ws.Cells.Item($r, 9).Value  = $vals.I   # I - Calidad
    $ws.Cells.Item($r, 10).Value = $vals.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Unidad de comercializacion
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Origen
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio $/Kg
}
